$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.359.60'
$ws.Range('E2').Value = '  -1.69%  '
$ws.Range('D3').Value = '1.592.65'
$ws.Range('E3').Value = '  -0.73%  '
$ws.Range('E4').Value = '  -0.42%  '
$ws.Range('D5').Value = "'210.11"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.92%  '
$ws.Range('D6').Value = "'0.505"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.53%  '
$ws.Range('E7').Value = '  -0.40%  '
$ws.Range('D8').Value = "'0.0611"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.44%  '
$ws.Range('E9').Value = '  -0.70%  '
$ws.Range('D10').Value = "'19.56"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.78%  '
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('D12').Value = '1.814.65'
$ws.Range('E12').Value = '  -0.84%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.600.27'
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = "'4.08"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.19%  '
$ws.Range('D15').Value = "'0.518"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.91%  '
$ws.Range('D16').Value = "'64.64"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.94%  '
$ws.Range('D17').Value = '26.364.30'
$ws.Range('E17').Value = '  -1.51%  '
$ws.Range('D18').Value = '0.0₃0728'
$ws.Range('E18').Value = '  -2.12%  '
$ws.Range('D19').Value = "'7.50"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.57%  '
$ws.Range('D20').Value = "'211.64"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('E21').Value = '  -0.39%  '
$ws.Range('E22').Value = '  -1.00%  '
$ws.Range('D23').Value = "'2.18"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.36%  '
$ws.Range('E24').Value = '  -1.66%  '
$ws.Range('D25').Value = "'145.17"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.67%  '
$ws.Range('E26').Value = '  -0.38%  '
$ws.Range('E27').Value = '  -1.49%  '
$ws.Range('E28').Value = '  -0.74%  '
$ws.Range('E29').Value = '  -0.71%  '
$ws.Range('E30').Value = '  -1.32%  '
$ws.Range('E31').Value = '  -0.40%  '
$ws.Range('E32').Value = '  -1.65%  '
$ws.Range('E33').Value = '  -0.36%  '
$ws.Range('D34').Value = '1.303.50'
$ws.Range('E34').Value = '  +1.52%  '
$ws.Range('D35').Value = "'0.615"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.10%  '
$ws.Range('E36').Value = '  -1.69%  '
$ws.Range('E37').Value = '  -1.28%  '
$ws.Range('E38').Value = '  -0.74%  '
$ws.Range('D39').Value = "'1.10"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -14.14%  '
$ws.Range('E40').Value = '  -2.06%  '
$ws.Range('E41').Value = '  -0.42%  '
$ws.Range('E42').Value = '  +2.77%  '
$ws.Range('D43').Value = "'62.71"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.35%  '
$ws.Range('D44').Value = "'2.14"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.86%  '
$ws.Range('E45').Value = '  -2.40%  '
$ws.Range('D46').Value = '1.728.54'
$ws.Range('E46').Value = '  -0.73%  '
$ws.Range('D47').Value = "'88.18"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.77%  '
$ws.Range('E48').Value = '  -4.56%  '
$ws.Range('E49').Value = '  +9.87%  '
$ws.Range('D50').Value = "'0.0985"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.36%  '
$ws.Range('E51').Value = '  -1.51%  '
